$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 7 (shifts old row7 -> row9)
$ws.Range("A7:T8").EntireRow.Insert()

# New row 7: Terminal Hortofrutícola Agro Chillán, Ñuble, 2023-12-15, Primera
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45275
$ws.Range("D7").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100101
$ws.Range("H7").Value = "Berries"
$ws.Range("I7").Value = 100101004
$ws.Range("J7").Value = "Frambuesa"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 9000
$ws.Range("P7").Value = 8500
$ws.Range("Q7").Value = "$/bandeja 2 kilos"
$ws.Range("R7").Value = "Región de Ñuble"
$ws.Range("S7").Value = 4250
$ws.Range("T7").Value = 2

# New row 8: Terminal Hortofrutícola Agro Chillán, Ñuble, 2023-12-15, Segunda
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 45275
$ws.Range("D8").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101004
$ws.Range("J8").Value = "Frambuesa"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Región de Ñuble"
$ws.Range("S8").Value = 3500
$ws.Range("T8").Value = 2
